$wb = $excel.ActiveWorkbook

# --- Sheet1: remove the now-unused "extra" rows (45-87) that only held
#     a leftover column-A counter series, shrinking the sheet back down
#     to the real data range (A1:N44). ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("45:87").Delete()

# --- Sheet1 becomes the active/selected sheet (was Sheet3 before), with
#     the view scrolled down near the bottom of the data and E64 selected. ---
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 1
$ws1.Range("E64").Select()

